$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new weekly price record as row 309 (shifting existing rows 309..356 down to 310..357),
# carrying the same fixed attributes as the rest of this market/category block plus the
# price/volume data from the week that used to sit at row 309, but with its own new date.
$ws.Rows.Item(309).Insert()

$ws.Cells.Item(309, 1).Value = 8
$ws.Cells.Item(309, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(309, 3).Value = "Coquimbo"
$ws.Cells.Item(309, 4).Value = 45154
$ws.Cells.Item(309, 5).Value = 4
$ws.Cells.Item(309, 6).Value = 100112037
$ws.Cells.Item(309, 7).Value = "Cebollín"
$ws.Cells.Item(309, 8).Value = "Sin especificar"
$ws.Cells.Item(309, 9).Value = "Primera"
$ws.Cells.Item(309, 10).Value = 1600
$ws.Cells.Item(309, 11).Value = 1000
$ws.Cells.Item(309, 12).Value = 1200
$ws.Cells.Item(309, 13).Value = 1100
$ws.Cells.Item(309, 14).Value = "`$/paquete 6 unidades"
$ws.Cells.Item(309, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(309, 16).Value = 183
$ws.Cells.Item(309, 17).Value = 6
$ws.Cells.Item(309, 18).Value = "Hortaliza"
